# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-07 01:18:33
#
# The "Recorded By" column (G) contains comma-separated lists of the
# users/systems that recorded a session. Upstream normalized the ordering
# of those names for a handful of rows (e.g. "dnasr281@gmail.com, System"
# became "System, dnasr281@gmail.com"). Re-apply the same exact text
# substitutions here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old "Recorded By" text -> new "Recorded By" text, exactly as seen
# in the upstream diff.
$replacements = @{
    'backup@backdoor.com, system, System' = 'backup@backdoor.com, System, system'
    'dnasr281@gmail.com, System'          = 'System, dnasr281@gmail.com'
    'dnasr281@gmail.com, admin@admin.com' = 'admin@admin.com, dnasr281@gmail.com'
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $value = $cell.Value2
    if ($null -ne $value -and $replacements.ContainsKey($value)) {
        $cell.Value2 = $replacements[$value]
    }
}
